$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1498.1
$ws.Range("I19").Value = 1699.3334
$ws.Range("K19").Value = 1699.3334
$ws.Range("M19").Value = -1524.3334

$ws.Range("H42").Value = 3402.6
$ws.Range("I42").Value = 415
$ws.Range("J42").Value = 5394.3335
$ws.Range("K42").Value = 1245
$ws.Range("L42").Value = 16183.0005
$ws.Range("M42").Value = -1015
$ws.Range("N42").Value = -16643.0005

$ws.Range("H98").Value = 1206.8572
$ws.Range("I98").Value = 887.4
$ws.Range("K98").Value = 887.4
$ws.Range("M98").Value = 610.6

$ws.Range("H112").Value = 4400.5
$ws.Range("J112").Value = 4489.4443
$ws.Range("L112").Value = 13468.3329
$ws.Range("N112").Value = -15684.3329

$ws.Range("H122").Value = 1206.8572
$ws.Range("I122").Value = 887.4
$ws.Range("K122").Value = 2662.2
$ws.Range("M122").Value = -212.1999999999998

$ws.Range("H125").Value = 24444
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 24444
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 219996
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -224916

$ws.Range("H135").Value = 4556.5
$ws.Range("I135").Value = 3312.9167
$ws.Range("J135").Value = 12018
$ws.Range("K135").Value = 29816.2503
$ws.Range("L135").Value = 108162
$ws.Range("M135").Value = -27281.2503
$ws.Range("N135").Value = -113232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3563.5715
$ws.Range("I122").Value = 3249
$ws.Range("J122").Value = 3983
$ws.Range("K122").Value = 9747
$ws.Range("L122").Value = 11949
$ws.Range("M122").Value = -7297
$ws.Range("N122").Value = -16849

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 1051
$ws.Range("I39").Value = 1051
$ws.Range("K39").Value = 1051
$ws.Range("M39").Value = -660

$ws.Range("H49").Value = 1051
$ws.Range("I49").Value = 1051
$ws.Range("K49").Value = 1051
$ws.Range("M49").Value = -869

$ws.Range("H107").Value = 1129
$ws.Range("I107").Value = 723.125
$ws.Range("K107").Value = 723.125
$ws.Range("M107").Value = 1196.875

$ws.Range("H122").Value = 1231.4286
$ws.Range("I122").Value = 1344
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 4032
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -1582
$ws.Range("N122").Value = -7750

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 3000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -3338

$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 3000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -3204

$ws.Range("H55").Value = 4199.6
$ws.Range("J55").Value = 4199.6
$ws.Range("L55").Value = 12598.8
$ws.Range("N55").Value = -12952.8

$ws.Range("H139").Value = 5000
$ws.Range("I139").Value = 5000
$ws.Range("K139").Value = 15000
$ws.Range("M139").Value = -9860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1506250
$ws.Range("I11").Value = 2004000
$ws.Range("K11").Value = 2004000
$ws.Range("M11").Value = -2003861

$ws.Range("H80").Value = 4499.6
$ws.Range("I80").Value = 3999.6667
$ws.Range("J80").Value = 5249.5
$ws.Range("K80").Value = 3999.6667
$ws.Range("L80").Value = 5249.5
$ws.Range("M80").Value = -3001.6667
$ws.Range("N80").Value = -7245.5

$ws.Range("H83").Value = 4499.6
$ws.Range("I83").Value = 3999.6667
$ws.Range("J83").Value = 5249.5
$ws.Range("K83").Value = 19998.3335
$ws.Range("L83").Value = 26247.5
$ws.Range("M83").Value = -15006.3335
$ws.Range("N83").Value = -36231.5

$ws.Range("H102").Value = 2164.8333
$ws.Range("I102").Value = 2164.8333
$ws.Range("K102").Value = 2164.8333
$ws.Range("M102").Value = -542.8332999999998

$ws.Range("H122").Value = 5044.1113
$ws.Range("I122").Value = 5132.6665
$ws.Range("J122").Value = 4999.8335
$ws.Range("K122").Value = 15397.9995
$ws.Range("L122").Value = 14999.5005
$ws.Range("M122").Value = -12947.9995
$ws.Range("N122").Value = -19899.5005

$ws.Range("H126").Value = 6400
$ws.Range("J126").Value = 7333.3335
$ws.Range("L126").Value = 22000.0005
$ws.Range("N126").Value = -26940.0005

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 3400
$ws.Range("J21").Value = 3400
$ws.Range("L21").Value = 3400
$ws.Range("N21").Value = -3748

$ws.Range("H22").Value = 1884.3125
$ws.Range("I22").Value = 2072.923
$ws.Range("J22").Value = 1067
$ws.Range("K22").Value = 2072.923
$ws.Range("L22").Value = 1067
$ws.Range("M22").Value = -1777.923
$ws.Range("N22").Value = -1657

$ws.Range("H27").Value = 1884.3125
$ws.Range("I27").Value = 2072.923
$ws.Range("J27").Value = 1067
$ws.Range("K27").Value = 2072.923
$ws.Range("L27").Value = 1067
$ws.Range("M27").Value = -1965.923
$ws.Range("N27").Value = -1281

$ws.Range("H40").Value = 4334.6665
$ws.Range("I40").Value = 3004
$ws.Range("K40").Value = 3004
$ws.Range("M40").Value = -2868

$ws.Range("H55").Value = 686.36365
$ws.Range("I55").Value = 656.9167
$ws.Range("J55").Value = 721.7
$ws.Range("K55").Value = 656.9167
$ws.Range("L55").Value = 721.7
$ws.Range("M55").Value = -483.9167
$ws.Range("N55").Value = -1067.7

$ws.Range("H68").Value = 2500
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751

$ws.Range("H71").Value = 2500
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756

$ws.Range("H122").Value = 6354.8
$ws.Range("I122").Value = 4561.5
$ws.Range("J122").Value = 7198.706
$ws.Range("K122").Value = 13684.5
$ws.Range("L122").Value = 21596.118
$ws.Range("M122").Value = -11234.5
$ws.Range("N122").Value = -26496.118

$ws.Range("H132").Value = 2828.8462
$ws.Range("I132").Value = 2472.75
$ws.Range("K132").Value = 7418.25
$ws.Range("M132").Value = -4888.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 49999
$ws.Range("J42").Value = 49999
$ws.Range("L42").Value = 49999
$ws.Range("N42").Value = -50755

$ws.Range("H70").Value = 34095
$ws.Range("I70").Value = 34095
$ws.Range("K70").Value = 34095
$ws.Range("M70").Value = -33780

$ws.Range("H73").Value = 34095
$ws.Range("I73").Value = 34095
$ws.Range("K73").Value = 34095
$ws.Range("M73").Value = -33003

$ws.Range("H136").Value = 909.5833
$ws.Range("I136").Value = 945.5
$ws.Range("J136").Value = 873.6667
$ws.Range("K136").Value = 2836.5
$ws.Range("L136").Value = 2621.0001
$ws.Range("M136").Value = -286.5
$ws.Range("N136").Value = -7721.0001
